$d = $word.ActiveDocument

# 1. Bump the document title version from 4.4 to 4.5
foreach ($p in $d.Paragraphs) {
    if ($p.Style.NameLocal -eq "Title") {
        $p.Range.Find.Execute("4.4", $true, $false, $false, $false, $false, $true, 1, $false, "4.5", 2)
        break
    }
}

# 2. Add a new Change Log row documenting the 4.5 release, right after the 4.4 row
$changeLogTable = $d.Tables.Item(1)
$newRow = $changeLogTable.Rows.Add()

$newRow.Cells.Item(1).Range.Text = "4.5"
$newRow.Cells.Item(2).Range.Text = "04/02/2021"
$newRow.Cells.Item(3).Range.Text = "Updated to release 4.5"

Write-Output "Edit complete"
